# Insert a new weekly price record for "Vega Modelo de Temuco - Chirimoya"
# at row 159, shifting the existing rows 159:183 down to 160:184.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 159 (pushes old row159..183 -> 160..184)
$ws.Rows.Item(159).Insert()

# Populate the new row 159 with the new weekly record
$ws.Cells.Item(159, 1).Value  = 10
$ws.Cells.Item(159, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(159, 3).Value  = "La Araucanía"
$ws.Cells.Item(159, 4).Value  = 45142
$ws.Cells.Item(159, 5).Value  = 9
$ws.Cells.Item(159, 6).Value  = "Fruta"
$ws.Cells.Item(159, 7).Value  = 100107
$ws.Cells.Item(159, 8).Value  = "Otros"
$ws.Cells.Item(159, 9).Value  = 100107002
$ws.Cells.Item(159, 10).Value = "Chirimoya"
$ws.Cells.Item(159, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(159, 12).Value = "Primera"
$ws.Cells.Item(159, 13).Value = 35
$ws.Cells.Item(159, 14).Value = 3500
$ws.Cells.Item(159, 15).Value = 3500
$ws.Cells.Item(159, 16).Value = 3500
$ws.Cells.Item(159, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(159, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(159, 19).Value = 3500
$ws.Cells.Item(159, 20).Value = 1
